$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$raw = @"
43918,1
43921,1
43925,6
43928,1
43935,3
43936,2
43937,7
43938,9
43939,8
43940,4
43941,7
43942,6
43943,8
43944,14
43945,8
43946,7
43947,4
43948,5
43949,23
43950,23
43951,16
43952,6
43953,6
43954,7
43955,8
43956,10
43957,7
43958,3
43959,6
43960,9
43962,4
43963,4
43964,7
43965,10
43966,15
43967,11
43969,22
43970,14
43971,13
43972,44
43973,35
43974,52
43975,2
43976,43
43977,14
43978,11
43979,34
43980,11
43981,22
43982,3
43983,14
43984,15
43985,14
43986,15
43987,22
43988,16
43989,3
43990,13
43991,16
43992,8
43993,62
43994,16
43995,9
43996,2
43997,27
43998,13
43999,22
44000,9
44001,23
44002,23
44003,1
44004,11
44005,13
44006,18
44007,28
44008,32
44009,18
44010,5
44011,19
44012,42
44013,22
44014,33
44015,27
44016,21
44017,3
44018,45
44019,47
44020,42
44021,76
44022,65
44023,70
44024,17
44025,56
44026,80
44027,62
44028,51
44029,93
44030,50
44031,38
44032,32
44033,88
44034,108
44035,65
44036,66
44037,95
44038,32
44039,137
44040,52
44041,100
44042,154
44043,114
44044,159
44045,56
44046,140
44047,186
44048,196
44049,201
44050,218
44051,172
44052,96
44053,127
44054,177
44055,258
44056,251
44057,278
44058,211
44059,57
44060,223
44061,174
44062,225
44063,403
44064,266
44065,259
44066,11
44067,293
44068,348
44069,197
44070,309
44071,24
44072,164
44073,26
44074,184
44075,107
44076,179
44077,53
44078,295
44079,103
44080,10
44081,282
44082,188
44083,132
44084,95
44085,263
44086,92
44087,15
44088,92
44089,96
44090,110
44091,131
44092,146
44093,38
44094,7
44095,62
44096,88
44097,165
44098,152
44099,143
44100,223
44101,8
44102,48
44103,39
44104,38
44105,31
44106,23
44107,22
44108,4
44109,24
44110,15
44111,28
44112,29
44113,45
44114,24
44115,5
44116,56
44117,49
44118,45
44119,31
44120,37
44121,16
44122,2
44123,42
44124,31
44125,31
44126,23
44127,32
44128,11
44129,1
44130,31
44131,53
44132,31
44133,17
44134,50
44135,13
44136,9
44137,35
44138,33
44139,40
44140,37
44142,46
44143,14
44144,56
44145,34
44146,43
44147,38
44148,39
44149,28
44150,4
44151,29
44152,37
44153,21
44154,32
44155,19
44156,20
44157,1
44158,37
44159,22
44160,47
44161,47
44162,46
44163,54
44164,4
44165,43
44166,23
44167,31
44168,17
44169,27
44170,14
44172,27
44173,10
44174,25
44175,26
44176,31
44177,16
44178,2
44179,43
44180,25
44181,17
44182,23
44183,15
44184,20
44185,2
44186,47
44187,23
44188,39
44189,9
44190,5
44191,26
44192,2
44193,41
44194,30
44195,58
44196,18
44198,14
44199,25
44200,2
44201,39
44202,26
44203,40
44204,73
44205,45
44206,7
44207,66
44208,75
44209,74
44210,80
44211,81
44212,78
44213,13
44214,121
44215,121
44216,102
44217,94
44218,107
44219,107
44220,17
44221,202
44222,138
44223,159
44224,130
44225,175
44226,138
44227,13
44228,185
44229,178
44230,147
44231,122
44232,126
44233,153
44234,17
44235,176
44236,188
44237,153
44238,129
44239,133
44240,109
44241,30
44242,130
44243,134
44244,109
44245,116
44246,103
44247,98
44248,22
44249,115
44250,93
44251,81
44252,124
44253,123
44254,79
44255,9
44256,116
44257,116
44258,99
44259,63
44260,63
44261,59
44262,13
44263,78
44264,78
44265,75
44266,68
44267,34
44268,32
44269,22
44270,78
44271,40
44272,87
44273,36
44274,49
44275,32
44276,5
44277,58
44278,46
44279,47
44280,55
44281,67
44282,38
44283,9
44284,64
44285,36
44286,49
44287,19
44288,21
44289,49
44290,16
44291,61
44292,61
44293,51
44294,54
44295,41
44296,42
44297,17
44298,49
44299,44
44300,59
44301,52
44302,72
44303,62
44304,21
44305,73
44306,52
44307,55
44308,71
44309,65
44310,65
44311,12
44312,71
44313,54
44314,33
44315,58
44316,68
44317,37
44318,3
44319,73
44320,50
44321,86
44322,61
44323,62
44324,48
44325,8
44326,97
44327,73
44328,72
44329,67
44330,50
44331,74
44332,17
44333,102
44334,84
44335,86
44336,87
44337,85
44338,60
44339,12
44340,103
44341,96
44342,26
44343,96
44344,91
44345,88
44346,19
44347,97
44348,94
44349,118
44350,98
44351,110
44352,64
44353,16
44354,87
44355,79
44356,68
44357,70
44358,62
44359,65
44360,10
44361,116
44362,138
44363,69
44364,80
44365,68
44366,72
44367,5
44368,122
44369,82
44370,72
44371,75
44372,78
44373,66
44374,31
44375,90
44376,26
44377,83
44378,72
44379,128
44380,39
44381,14
44382,83
44383,66
44384,52
44385,63
44386,59
44387,52
44388,18
44389,66
44390,59
44391,55
44392,37
44393,51
44394,44
44395,23
44396,53
44397,46
44398,43
44399,41
44400,28
44401,27
44402,7
44403,44
44404,24
44405,19
44406,11
44407,38
44408,37
44409,14
44410,36
44411,30
44412,26
44413,28
44414,27
44415,18
44416,10
44417,20
44418,13
44419,17
44420,19
44421,18
44422,15
44423,5
44424,25
44425,8
44426,15
44427,17
44428,11
44429,11
44430,2
44431,9
44432,20
44433,7
44434,9
44435,3
44436,4
44437,3
44438,8
44439,4
44440,3
44441,8
44442,8
44443,2
44445,11
44446,8
44447,7
44448,2
44450,4
44451,2
44452,2
44453,6
44454,6
44455,2
44456,6
44457,8
44458,3
44459,7
44460,4
44461,6
44462,4
44463,4
44464,5
44466,12
44467,6
44468,9
44469,13
44470,9
44471,7
44472,1
44473,6
44474,7
44475,6
44476,11
44477,1
44478,11
44479,1
44480,11
44481,18
44482,13
44483,10
44484,5
44485,12
44486,2
44487,16
44488,19
44489,11
44490,5
44491,20
44492,17
44493,3
44494,21
44495,12
44496,9
44497,25
44498,30
44499,14
44500,2
44501,1
44502,33
44503,31
44504,22
44505,25
44506,29
44507,4
44508,31
44509,52
44510,39
44511,29
44512,20
44513,26
44514,4
44515,27
44516,24
44517,31
44518,33
44519,28
44520,14
44521,4
44522,17
44523,29
44524,37
44525,28
44526,17
44527,18
44528,1
44529,23
44530,32
44531,13
44532,41
44533,46
44534,23
44535,5
44536,26
44537,34
44538,2
44539,29
44540,28
44541,25
44542,8
44543,42
44544,33
44545,20
44546,51
44547,52
44548,24
44549,3
44550,38
44551,43
44552,21
44553,14
44554,48
44555,10
44556,5
44557,44
44558,52
44559,54
44560,53
44561,74
44562,11
44563,14
44564,86
44565,129
44566,177
44567,220
44568,326
44569,325
44570,88
44571,410
44572,422
44573,673
44574,719
44575,772
44576,748
44577,220
44578,1126
44579,1087
44580,1228
44581,1040
44582,1119
44583,853
44584,233
44585,1124
44586,874
44587,693
44588,816
44589,695
44590,431
44591,293
44592,653
44593,470
44594,502
44595,395
44596,344
44597,279
44598,72
44599,340
44600,242
44601,206
44602,184
44603,173
44604,110
44605,16
44606,119
44607,107
44608,49
"@

$lines = $raw -split "`n"
$n = $lines.Length
Write-Output "rows=$n"

$data = New-Object 'object[,]' $n,2
for ($i = 0; $i -lt $n; $i++) {
    $parts = $lines[$i].Trim() -split ","
    $data[$i,0] = [double]$parts[0]
    $data[$i,1] = [double]$parts[1]
}

$lastRow = $n + 1
$rng = $ws.Range("A2:B$lastRow")
$rng.Value2 = $data

# Ensure the new rows beyond the old 627 (rows 628..671) get the same
# date-column style (style "2": custom date format, bold, centered, bordered)
# as the rest of column A, by copying format down from an already-styled cell.
$oldLastRow = 627
if ($lastRow -gt $oldLastRow) {
    $srcFormat = $ws.Range("A$oldLastRow")
    $dstFormat = $ws.Range("A" + ($oldLastRow + 1) + ":A" + $lastRow)
    $srcFormat.Copy()
    $dstFormat.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

Write-Output "Done. Used range: $($ws.UsedRange.Address())"
